# Update countries & provincias Spain
#
# Refreshes the COVID-19 "Pais" snapshot: the timestamp banner, several
# countries' case statistics, and the ranking order for a handful of
# countries whose total-case counts crossed each other (so their rows
# swapped places in the leaderboard).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 12:50"

# Country list was re-sorted (ranking refresh): update the displayed names
# for the rows whose country changed position.
$ws.Range("A86").Value = "Bielorrusia"  # was Taiwan
$ws.Range("A87").Value = "Taiwan"  # was Reunion
$ws.Range("A88").Value = "Reunion"  # was Camerun
$ws.Range("A89").Value = "Camerun"  # was Bielorrusia
$ws.Range("A100").Value = "Senegal"  # was Ghana
$ws.Range("A101").Value = "Ghana"  # was Malta
$ws.Range("A102").Value = "Malta"  # was Senegal
$ws.Range("A151").Value = "San Martin (Parte Holandesa)"  # was Eritrea
$ws.Range("A152").Value = "Eritrea"  # was San Martin (Parte Francesa)
$ws.Range("A153").Value = "San Martin (Parte Francesa)"  # was Congo
$ws.Range("A154").Value = "Congo"  # was Gabon
$ws.Range("A155").Value = "Gabon"  # was Birmania
$ws.Range("A156").Value = "Birmania"  # was Tanzania
$ws.Range("A157").Value = "Tanzania"  # was Guyana
$ws.Range("A158").Value = "Guyana"  # was Maldivas
$ws.Range("A159").Value = "Maldivas"  # was Nueva Caledonia
$ws.Range("A160").Value = "Nueva Caledonia"  # was Haiti
$ws.Range("A161").Value = "Haiti"  # was San Martin (Parte Holandesa)
$ws.Range("A184").Value = "Republica del Chad"  # was Republica de Africa Central
$ws.Range("A185").Value = "Republica de Africa Central"  # was Republica del Chad
$ws.Range("A188").Value = "Santa Sede"  # was Fiyi
$ws.Range("A189").Value = "Fiyi"  # was Santa Sede
$ws.Range("A190").Value = "Liberia"  # was Nepal
$ws.Range("A191").Value = "Cabo Verde"  # was Liberia
$ws.Range("A192").Value = "Nepal"  # was Cabo Verde

# Updated case statistics for the affected rows
$ws.Range("E21").Value = 6655  # was 6656
$ws.Range("G21").Value = 1  # was 0
$ws.Range("H21").Value = 37  # was 36
$ws.Range("D31").Value = 283  # was 267
$ws.Range("E31").Value = 2784  # was 2800
$ws.Range("F31").Value = 83  # was 78
$ws.Range("D35").Value = 514  # was 472
$ws.Range("E35").Value = 2040  # was 2082
$ws.Range("D66").Value = 49  # was 34
$ws.Range("E66").Value = 639  # was 654
$ws.Range("B69").Value = 672  # was 643
$ws.Range("C69").Value = 29  # was 0
$ws.Range("D69").Value = 382  # was 381
$ws.Range("E69").Value = 286  # was 258
$ws.Range("B86").Value = 351  # was 348
$ws.Range("C86").Value = 47  # was 9
$ws.Range("D86").Value = 53  # was 50
$ws.Range("E86").Value = 294  # was 293
$ws.Range("F86").Value = 11  # was 0
$ws.Range("H86").Value = 4  # was 5
$ws.Range("B87").Value = 348  # was 308
$ws.Range("C87").Value = 9  # was 0
$ws.Range("D87").Value = 50  # was 40
$ws.Range("E87").Value = 293  # was 268
$ws.Range("F87").Value = 0  # was 3
$ws.Range("H87").Value = 5  # was 0
$ws.Range("B88").Value = 308  # was 306
$ws.Range("D88").Value = 40  # was 10
$ws.Range("E88").Value = 268  # was 289
$ws.Range("F88").Value = 3  # was 0
$ws.Range("H88").Value = 0  # was 7
$ws.Range("B89").Value = 306  # was 304
$ws.Range("D89").Value = 10  # was 53
$ws.Range("E89").Value = 289  # was 247
$ws.Range("F89").Value = 0  # was 11
$ws.Range("H89").Value = 7  # was 4
$ws.Range("B100").Value = 207  # was 204
$ws.Range("C100").Value = 12  # was 0
$ws.Range("D100").Value = 66  # was 31
$ws.Range("E100").Value = 140  # was 168
$ws.Range("F100").Value = 1  # was 2
$ws.Range("H100").Value = 1  # was 5
$ws.Range("B101").Value = 204  # was 196
$ws.Range("D101").Value = 31  # was 2
$ws.Range("E101").Value = 168  # was 194
$ws.Range("H101").Value = 5  # was 0
$ws.Range("B102").Value = 202  # was 195
$ws.Range("C102").Value = 6  # was 0
$ws.Range("D102").Value = 2  # was 55
$ws.Range("E102").Value = 200  # was 139
$ws.Range("F102").Value = 2  # was 0
$ws.Range("H102").Value = 0  # was 1
$ws.Range("B104").Value = 190  # was 184
$ws.Range("C104").Value = 6  # was 0
$ws.Range("E104").Value = 168  # was 162
$ws.Range("B151").Value = 23  # was 22
$ws.Range("C151").Value = 5  # was 0
$ws.Range("D151").Value = 6  # was 0
$ws.Range("E151").Value = 15  # was 22
$ws.Range("G151").Value = 1  # was 0
$ws.Range("H151").Value = 2  # was 0
$ws.Range("D152").Value = 0  # was 2
$ws.Range("E152").Value = 22  # was 19
$ws.Range("H152").Value = 0  # was 1
$ws.Range("E153").Value = 19  # was 18
$ws.Range("H153").Value = 1  # was 2
$ws.Range("B154").Value = 22  # was 21
$ws.Range("D154").Value = 2  # was 1
$ws.Range("E154").Value = 18  # was 19
$ws.Range("H154").Value = 2  # was 1
$ws.Range("B155").Value = 21  # was 20
$ws.Range("D155").Value = 1  # was 0
$ws.Range("D156").Value = 0  # was 2
$ws.Range("E156").Value = 19  # was 17
$ws.Range("B157").Value = 20  # was 19
$ws.Range("D157").Value = 3  # was 0
$ws.Range("E157").Value = 16  # was 15
$ws.Range("H157").Value = 1  # was 4
$ws.Range("D158").Value = 0  # was 13
$ws.Range("E158").Value = 15  # was 6
$ws.Range("H158").Value = 4  # was 0
$ws.Range("B159").Value = 19  # was 18
$ws.Range("D159").Value = 13  # was 1
$ws.Range("E159").Value = 6  # was 17
$ws.Range("C160").Value = 0  # was 2
$ws.Range("C161").Value = 2  # was 0
$ws.Range("D161").Value = 1  # was 6
$ws.Range("E161").Value = 17  # was 11
$ws.Range("H161").Value = 0  # was 1
$ws.Range("C184").Value = 0  # was 5
$ws.Range("C185").Value = 5  # was 0
$ws.Range("B190").Value = 6  # was 7
$ws.Range("C190").Value = 0  # was 1
$ws.Range("D190").Value = 0  # was 1
$ws.Range("E191").Value = 5  # was 6
$ws.Range("H191").Value = 1  # was 0
$ws.Range("D192").Value = 1  # was 0
$ws.Range("H192").Value = 0  # was 1
